$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 48, shifting the existing rows 48:110 down to 49:111
$ws.Rows.Item(48).Insert()

# Populate the newly inserted row 48 with the new weekly data record
$ws.Cells.Item(48, 1).Value = 3
$ws.Cells.Item(48, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(48, 3).Value = "Coquimbo"
$ws.Cells.Item(48, 4).Value = 44546
$ws.Cells.Item(48, 5).Value = 5
$ws.Cells.Item(48, 6).Value = 100112030
$ws.Cells.Item(48, 7).Value = "Poroto granado"
$ws.Cells.Item(48, 8).Value = "Sin especificar"
$ws.Cells.Item(48, 9).Value = "Primera"
$ws.Cells.Item(48, 10).Value = 73
$ws.Cells.Item(48, 11).Value = 44000
$ws.Cells.Item(48, 12).Value = 45000
$ws.Cells.Item(48, 13).Value = 44521
$ws.Cells.Item(48, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(48, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(48, 16).Value = 1781
$ws.Cells.Item(48, 17).Value = 25
$ws.Cells.Item(48, 18).Value = "Hortaliza"
